$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The PCB now uses 2x ADXRS620 gyros instead of 1
$ws.Cells.Item(5, 5).Value = 2

# Make room above "DIGITAL BOARD 3.3V Current" for the new 1.8V analog rail
# summary row (shifts old rows 23/24 down to 24/25, leaving a gap at 23).
$ws.Rows.Item(23).Insert()

# New summary row: ANALOG BOARD 1.8V current
$ws.Cells.Item(20, 1).Value = "ANALOG BOARD 1.8V current"
$ws.Cells.Item(20, 2).Value = 5

# The 1.8V regulator (TPS-1.8, row 6 / F6) now has its own summary line above,
# so it should no longer be folded into the 3.3V analog total.
$ws.Cells.Item(21, 2).Formula = "=SUM(F2,F8)"

# Document the 1.8V regulator part + its power-dissipation note next to the
# ANALOG BOARD 5V Current row.
$ws.Cells.Item(22, 4).Value = "MCP1791T-5002E/DCCT-ND"
$ws.Cells.Item(22, 5).Value = "Max current 70mA, max power at 120 degrees F is 1.22 Watts, we are dissipating < .6W"

# Widen the columns to fit the new part number / note text.
$ws.Columns.Item(4).ColumnWidth = 25
$ws.Columns.Item(5).ColumnWidth = 42

$excel.Calculate()

$ws.Range("E23").Select() | Out-Null
